$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1.2
$ws.Range("I11").Value = 1.2
$ws.Range("K11").Value = 1.2
$ws.Range("M11").Value = 138.8

$ws.Range("H76").Value = 3987.5
$ws.Range("I76").Value = 3987.5
$ws.Range("K76").Value = 3987.5
$ws.Range("M76").Value = -3672.5

$ws.Range("H79").Value = 3987.5
$ws.Range("I79").Value = 3987.5
$ws.Range("K79").Value = 3987.5
$ws.Range("M79").Value = -2895.5

$ws.Range("H80").Value = 292.16666
$ws.Range("I80").Value = 157.90909
$ws.Range("J80").Value = 503.14285
$ws.Range("K80").Value = 473.72727
$ws.Range("L80").Value = 1509.42855
$ws.Range("M80").Value = 524.27273
$ws.Range("N80").Value = -3505.42855

$ws.Range("H83").Value = 292.16666
$ws.Range("I83").Value = 157.90909
$ws.Range("J83").Value = 503.14285
$ws.Range("K83").Value = 1421.18181
$ws.Range("L83").Value = 4528.28565
$ws.Range("M83").Value = 3570.81819
$ws.Range("N83").Value = -14512.28565

$ws.Range("H98").Value = 1685
$ws.Range("I98").Value = 1660.7693
$ws.Range("K98").Value = 1660.7693
$ws.Range("M98").Value = -162.7692999999999

$ws.Range("H122").Value = 1685
$ws.Range("I122").Value = 1660.7693
$ws.Range("K122").Value = 4982.3079
$ws.Range("M122").Value = -2532.3079

$ws.Range("H132").Value = 39485.445
$ws.Range("I132").Value = 42090.965
$ws.Range("K132").Value = 126272.895
$ws.Range("M132").Value = -123742.895

$ws.Range("H135").Value = 1280.7858
$ws.Range("I135").Value = 1323.7778
$ws.Range("K135").Value = 11914.0002
$ws.Range("M135").Value = -9379.0002

$ws.Range("H137").Value = 1516373
$ws.Range("J137").Value = 1648555
$ws.Range("L137").Value = 4945665
$ws.Range("N137").Value = -4950765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1350.5454
$ws.Range("I2").Value = 1249.5111
$ws.Range("J2").Value = 1805.2
$ws.Range("K2").Value = 1249.5111
$ws.Range("L2").Value = 1805.2
$ws.Range("M2").Value = -1136.5111
$ws.Range("N2").Value = -2031.2

$ws.Range("H32").Value = 1867.87
$ws.Range("I32").Value = 1836.2323
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1836.2323
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -1549.2323
$ws.Range("N32").Value = -5574

$ws.Range("H45").Value = 4308.174
$ws.Range("I45").Value = 3742.2856
$ws.Range("J45").Value = 10250
$ws.Range("K45").Value = 3742.2856
$ws.Range("L45").Value = 10250
$ws.Range("M45").Value = -3365.2856
$ws.Range("N45").Value = -11004

$ws.Range("H63").Value = 3543.8
$ws.Range("I63").Value = 2014.5454
$ws.Range("K63").Value = 2014.5454
$ws.Range("M63").Value = -1328.5454

$ws.Range("H66").Value = 3543.8
$ws.Range("I66").Value = 2014.5454
$ws.Range("K66").Value = 10072.727
$ws.Range("M66").Value = -6640.726999999999

$ws.Range("H116").Value = 1350.5454
$ws.Range("I116").Value = 1249.5111
$ws.Range("J116").Value = 1805.2
$ws.Range("K116").Value = 1249.5111
$ws.Range("L116").Value = 1805.2
$ws.Range("M116").Value = 1044.4889
$ws.Range("N116").Value = -6393.2

$ws.Range("H132").Value = 271353.03
$ws.Range("I132").Value = 433040.84
$ws.Range("K132").Value = 1299122.52
$ws.Range("M132").Value = -1296592.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1350.5454
$ws.Range("I3").Value = 1249.5111
$ws.Range("J3").Value = 1805.2
$ws.Range("K3").Value = 1249.5111
$ws.Range("L3").Value = 1805.2
$ws.Range("M3").Value = -1135.5111
$ws.Range("N3").Value = -2033.2

$ws.Range("H20").Value = 1546.8
$ws.Range("I20").Value = 1420.875
$ws.Range("K20").Value = 1420.875
$ws.Range("M20").Value = -1173.875

$ws.Range("H94").Value = 563.04
$ws.Range("I94").Value = 544.8333
$ws.Range("K94").Value = 544.8333
$ws.Range("M94").Value = -93.83330000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7401.702
$ws.Range("I31").Value = 1752.0526
$ws.Range("K31").Value = 1752.0526
$ws.Range("M31").Value = -1457.0526

$ws.Range("H34").Value = 7401.702
$ws.Range("I34").Value = 1752.0526
$ws.Range("K34").Value = 1752.0526
$ws.Range("M34").Value = -1550.0526

$ws.Range("H56").Value = 25000
$ws.Range("J56").Value = 25000
$ws.Range("L56").Value = 25000
$ws.Range("N56").Value = -26690

$ws.Range("H58").Value = 414000.6
$ws.Range("I58").Value = 538783.4
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 538783.4
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -538580.4
$ws.Range("N58").Value = -4406

$ws.Range("H132").Value = 10886993
$ws.Range("I132").Value = 21452.334
$ws.Range("K132").Value = 64357.00199999999
$ws.Range("M132").Value = -61827.00199999999

$ws.Range("H136").Value = 414000.6
$ws.Range("I136").Value = 538783.4
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 1616350.2
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1613800.2
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 284.42856
$ws.Range("I44").Value = 446
$ws.Range("J44").Value = 163.25
$ws.Range("K44").Value = 1338
$ws.Range("L44").Value = 489.75
$ws.Range("M44").Value = -940
$ws.Range("N44").Value = -1285.75

$ws.Range("H106").Value = 6135.5557
$ws.Range("I106").Value = 6574.3335
$ws.Range("K106").Value = 19723.0005
$ws.Range("M106").Value = -18777.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 209
$ws.Range("I2").Value = 110.25
$ws.Range("K2").Value = 110.25
$ws.Range("M2").Value = 2.75

$ws.Range("H102").Value = 2407.9429
$ws.Range("I102").Value = 2053.8572
$ws.Range("K102").Value = 2053.8572
$ws.Range("M102").Value = -431.8571999999999

$ws.Range("H113").Value = 1756
$ws.Range("I113").Value = 1416.875
$ws.Range("J113").Value = 3112.5
$ws.Range("K113").Value = 1416.875
$ws.Range("L113").Value = 3112.5
$ws.Range("M113").Value = 753.125
$ws.Range("N113").Value = -7452.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7984.5454
$ws.Range("I7").Value = 7783
$ws.Range("K7").Value = 7783
$ws.Range("M7").Value = -7671

$ws.Range("H16").Value = 50170
$ws.Range("I16").Value = 50170
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 50170
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -50000
$ws.Range("N16").ClearContents()

$ws.Range("H92").Value = 60389
$ws.Range("J92").Value = 60389
$ws.Range("L92").Value = 60389
$ws.Range("N92").Value = -65381

$ws.Range("H122").Value = 33788.12
$ws.Range("I122").Value = 2950.3044
$ws.Range("K122").Value = 8850.913199999999
$ws.Range("M122").Value = -6400.913199999999

$ws.Range("H126").Value = 7984.5454
$ws.Range("I126").Value = 7783
$ws.Range("K126").Value = 23349
$ws.Range("M126").Value = -20879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1233
$ws.Range("I81").Value = 1233
$ws.Range("K81").Value = 2466
$ws.Range("M81").Value = -1405

$ws.Range("H84").Value = 1233
$ws.Range("I84").Value = 1233
$ws.Range("K84").Value = 12330
$ws.Range("M84").Value = -7026

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H113").Value = 1083.3636
$ws.Range("I113").Value = 895.0909
$ws.Range("K113").Value = 2685.2727
$ws.Range("M113").Value = -515.2727

$ws.Range("H132").Value = 6496386
$ws.Range("I132").Value = 14379891
$ws.Range("J132").Value = 4087.4119
$ws.Range("K132").Value = 43139673
$ws.Range("L132").Value = 12262.2357
$ws.Range("M132").Value = -43137143
$ws.Range("N132").Value = -17322.2357

